$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets("ALC")
$ws_ALC.Range("H53").Value = 142.91667
$ws_ALC.Range("I53").Value = 176.875
$ws_ALC.Range("J53").Value = 75
$ws_ALC.Range("K53").Value = 176.875
$ws_ALC.Range("L53").Value = 75
$ws_ALC.Range("M53").Value = 460.125
$ws_ALC.Range("N53").Value = -1349
$ws_ALC.Range("H132").Value = 1544313
$ws_ALC.Range("I132").Value = 1736842.8
$ws_ALC.Range("J132").Value = 4074.75
$ws_ALC.Range("K132").Value = 5210528.4
$ws_ALC.Range("L132").Value = 12224.25
$ws_ALC.Range("M132").Value = -5207998.4
$ws_ALC.Range("N132").Value = -17284.25
$ws_ALC.Range("H137").Value = 1716
$ws_ALC.Range("I137").Value = 1633.3334
$ws_ALC.Range("K137").Value = 4900.0002
$ws_ALC.Range("M137").Value = -2350.0002
$ws_ALC.Range("H138").Value = 2352.5205
$ws_ALC.Range("J138").Value = 4509.1665
$ws_ALC.Range("L138").Value = 13527.4995
$ws_ALC.Range("N138").Value = -23807.4995
$ws_ALC.Range("H141").Value = 4690.1113
$ws_ALC.Range("I141").Value = 2891.1765
$ws_ALC.Range("J141").Value = 5782.3213
$ws_ALC.Range("K141").Value = 8673.529500000001
$ws_ALC.Range("L141").Value = 17346.9639
$ws_ALC.Range("M141").Value = -3493.529500000001
$ws_ALC.Range("N141").Value = -27706.9639
$ws_ARM = $wb.Worksheets("ARM")
$ws_ARM.Range("H32").Value = 9197.471
$ws_ARM.Range("I32").Value = 5214.075
$ws_ARM.Range("K32").Value = 5214.075
$ws_ARM.Range("M32").Value = -4927.075
$ws_ARM.Range("H45").Value = 1886.5385
$ws_ARM.Range("I45").Value = 1995.5
$ws_ARM.Range("J45").Value = 1523.3334
$ws_ARM.Range("K45").Value = 1995.5
$ws_ARM.Range("L45").Value = 1523.3334
$ws_ARM.Range("M45").Value = -1618.5
$ws_ARM.Range("N45").Value = -2277.3334
$ws_ARM.Range("H61").Value = 1571.9333
$ws_ARM.Range("I61").Value = 1350.409
$ws_ARM.Range("J61").Value = 2181.125
$ws_ARM.Range("K61").Value = 1350.409
$ws_ARM.Range("L61").Value = 2181.125
$ws_ARM.Range("M61").Value = -1138.409
$ws_ARM.Range("N61").Value = -2605.125
$ws_ARM.Range("H74").Value = 604.5536
$ws_ARM.Range("I74").Value = 525.4375
$ws_ARM.Range("J74").Value = 1079.25
$ws_ARM.Range("K74").Value = 525.4375
$ws_ARM.Range("L74").Value = 1079.25
$ws_ARM.Range("M74").Value = 348.5625
$ws_ARM.Range("N74").Value = -2827.25
$ws_ARM.Range("H77").Value = 604.5536
$ws_ARM.Range("I77").Value = 525.4375
$ws_ARM.Range("J77").Value = 1079.25
$ws_ARM.Range("K77").Value = 2627.1875
$ws_ARM.Range("L77").Value = 5396.25
$ws_ARM.Range("M77").Value = 1740.8125
$ws_ARM.Range("N77").Value = -14132.25
$ws_ARM.Range("H132").Value = 1947.0938
$ws_ARM.Range("I132").Value = 1411.0588
$ws_ARM.Range("J132").Value = 2554.6
$ws_ARM.Range("K132").Value = 4233.1764
$ws_ARM.Range("L132").Value = 7663.799999999999
$ws_ARM.Range("M132").Value = -1703.1764
$ws_ARM.Range("N132").Value = -12723.8
$ws_ARM.Range("H136").Value = 1571.9333
$ws_ARM.Range("I136").Value = 1350.409
$ws_ARM.Range("J136").Value = 2181.125
$ws_ARM.Range("K136").Value = 4051.227
$ws_ARM.Range("L136").Value = 6543.375
$ws_ARM.Range("M136").Value = -1501.227
$ws_ARM.Range("N136").Value = -11643.375
$ws_BSM = $wb.Worksheets("BSM")
$ws_BSM.Range("H20").Value = 29421736
$ws_BSM.Range("I20").Value = 35724836
$ws_BSM.Range("J20").Value = 7266.6665
$ws_BSM.Range("K20").Value = 35724836
$ws_BSM.Range("L20").Value = 7266.6665
$ws_BSM.Range("M20").Value = -35724589
$ws_BSM.Range("N20").Value = -7760.6665
$ws_BSM.Range("H134").Value = 2914.2856
$ws_BSM.Range("I134").Value = 2880
$ws_BSM.Range("J134").Value = 3000
$ws_BSM.Range("K134").Value = 8640
$ws_BSM.Range("L134").Value = 9000
$ws_BSM.Range("M134").Value = -6105
$ws_BSM.Range("N134").Value = -14070
$ws_CRP = $wb.Worksheets("CRP")
$ws_CRP.Range("H31").Value = 2759.4075
$ws_CRP.Range("I31").Value = 2492.2856
$ws_CRP.Range("K31").Value = 2492.2856
$ws_CRP.Range("M31").Value = -2197.2856
$ws_CRP.Range("H34").Value = 2759.4075
$ws_CRP.Range("I34").Value = 2492.2856
$ws_CRP.Range("K34").Value = 2492.2856
$ws_CRP.Range("M34").Value = -2290.2856
$ws_CRP.Range("H58").Value = 1036.3438
$ws_CRP.Range("I58").Value = 706.5
$ws_CRP.Range("J58").Value = 1586.0834
$ws_CRP.Range("K58").Value = 706.5
$ws_CRP.Range("L58").Value = 1586.0834
$ws_CRP.Range("M58").Value = -503.5
$ws_CRP.Range("N58").Value = -1992.0834
$ws_CRP.Range("H102").Value = 25115
$ws_CRP.Range("I102").Value = 20219
$ws_CRP.Range("J102").Value = 26747
$ws_CRP.Range("K102").Value = 20219
$ws_CRP.Range("L102").Value = 26747
$ws_CRP.Range("M102").Value = -17785
$ws_CRP.Range("N102").Value = -31615
$ws_CRP.Range("H132").Value = 6009.5186
$ws_CRP.Range("I132").Value = 6539.95
$ws_CRP.Range("K132").Value = 19619.85
$ws_CRP.Range("M132").Value = -17089.85
$ws_CRP.Range("H136").Value = 1036.3438
$ws_CRP.Range("I136").Value = 706.5
$ws_CRP.Range("J136").Value = 1586.0834
$ws_CRP.Range("K136").Value = 2119.5
$ws_CRP.Range("L136").Value = 4758.2502
$ws_CRP.Range("M136").Value = 430.5
$ws_CRP.Range("N136").Value = -9858.2502
$ws_CUL = $wb.Worksheets("CUL")
$ws_CUL.Range("H113").Value = 407.34692
$ws_CUL.Range("J113").Value = 420.8
$ws_CUL.Range("L113").Value = 1262.4
$ws_CUL.Range("N113").Value = -5602.4
$ws_CUL.Range("H122").Value = 6539.3887
$ws_CUL.Range("I122").Value = 595
$ws_CUL.Range("J122").Value = 8237.786
$ws_CUL.Range("K122").Value = 5355
$ws_CUL.Range("L122").Value = 74140.07399999999
$ws_CUL.Range("M122").Value = -2905
$ws_CUL.Range("N122").Value = -79040.07399999999
$ws_CUL.Range("H131").Value = 930
$ws_CUL.Range("I131").Value = 431.23077
$ws_CUL.Range("J131").Value = 1020.05554
$ws_CUL.Range("K131").Value = 1293.69231
$ws_CUL.Range("L131").Value = 3060.16662
$ws_CUL.Range("M131").Value = 3746.30769
$ws_CUL.Range("N131").Value = -13140.16662
$ws_CUL.Range("H141").Value = 4286
$ws_CUL.Range("I141").Value = 1562.5
$ws_CUL.Range("J141").Value = 9733
$ws_CUL.Range("K141").Value = 4687.5
$ws_CUL.Range("L141").Value = 29199
$ws_CUL.Range("M141").Value = 492.5
$ws_CUL.Range("N141").Value = -39559
$ws_GSM = $wb.Worksheets("GSM")
$ws_GSM.Range("H15").Value = 10330.25
$ws_GSM.Range("J15").Value = 10330.25
$ws_GSM.Range("L15").Value = 10330.25
$ws_GSM.Range("N15").Value = -10906.25
$ws_GSM.Range("H70").Value = 5551.6
$ws_GSM.Range("I70").Value = 5122.385
$ws_GSM.Range("K70").Value = 5122.385
$ws_GSM.Range("M70").Value = -4852.385
$ws_GSM.Range("H73").Value = 5551.6
$ws_GSM.Range("I73").Value = 5122.385
$ws_GSM.Range("K73").Value = 5122.385
$ws_GSM.Range("M73").Value = -4186.385
$ws_GSM.Range("H80").Value = 2350
$ws_GSM.Range("J80").Value = 2350
$ws_GSM.Range("L80").Value = 2350
$ws_GSM.Range("N80").Value = -4346
$ws_GSM.Range("H81").Value = 10330.25
$ws_GSM.Range("J81").Value = 10330.25
$ws_GSM.Range("L81").Value = 10330.25
$ws_GSM.Range("N81").Value = -12326.25
$ws_GSM.Range("H83").Value = 2350
$ws_GSM.Range("J83").Value = 2350
$ws_GSM.Range("L83").Value = 11750
$ws_GSM.Range("N83").Value = -21734
$ws_GSM.Range("H84").Value = 10330.25
$ws_GSM.Range("J84").Value = 10330.25
$ws_GSM.Range("L84").Value = 30990.75
$ws_GSM.Range("N84").Value = -40974.75
$ws_GSM.Range("H122").Value = 4104
$ws_GSM.Range("I122").Value = 0
$ws_GSM.Range("J122").Value = 4104
$ws_GSM.Range("K122").Value = 0
$ws_GSM.Range("L122").Value = 12312
$ws_GSM.Range("M122").ClearContents()
$ws_GSM.Range("N122").Value = -17212
$ws_LTW = $wb.Worksheets("LTW")
$ws_LTW.Range("H40").Value = 1990.0667
$ws_LTW.Range("I40").Value = 1982.1282
$ws_LTW.Range("J40").Value = 2041.6666
$ws_LTW.Range("K40").Value = 1982.1282
$ws_LTW.Range("L40").Value = 2041.6666
$ws_LTW.Range("M40").Value = -1846.1282
$ws_LTW.Range("N40").Value = -2313.6666
$ws_LTW.Range("H99").Value = 28000
$ws_LTW.Range("I99").Value = 0
$ws_LTW.Range("J99").Value = 28000
$ws_LTW.Range("K99").Value = 0
$ws_LTW.Range("L99").Value = 28000
$ws_LTW.Range("M99").ClearContents()
$ws_LTW.Range("N99").Value = -33990
$ws_LTW.Range("H122").Value = 3020.1025
$ws_LTW.Range("I122").Value = 2930.75
$ws_LTW.Range("J122").Value = 3428.5715
$ws_LTW.Range("K122").Value = 8792.25
$ws_LTW.Range("L122").Value = 10285.7145
$ws_LTW.Range("M122").Value = -6342.25
$ws_LTW.Range("N122").Value = -15185.7145
$ws_LTW.Range("H136").Value = 1994.6538
$ws_LTW.Range("I136").Value = 834.75
$ws_LTW.Range("J136").Value = 2510.1667
$ws_LTW.Range("K136").Value = 2504.25
$ws_LTW.Range("L136").Value = 7530.500100000001
$ws_LTW.Range("M136").Value = 45.75
$ws_LTW.Range("N136").Value = -12630.5001
$ws_WVR = $wb.Worksheets("WVR")
$ws_WVR.Range("H102").Value = 50000
$ws_WVR.Range("J102").Value = 50000
$ws_WVR.Range("L102").Value = 50000
$ws_WVR.Range("N102").Value = -56490
$ws_WVR.Range("H136").Value = 4274.6772
$ws_WVR.Range("I136").Value = 843.8570999999999
$ws_WVR.Range("J136").Value = 7100.0586
$ws_WVR.Range("K136").Value = 2531.5713
$ws_WVR.Range("L136").Value = 21300.1758
$ws_WVR.Range("M136").Value = 18.42870000000039
$ws_WVR.Range("N136").Value = -26400.1758
